# Auto-generated edit script updating cryptocurrency Price (D) and
# Volume(1h) (E) columns to match the refreshed scrape from the
# "Updated cryptos list" GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.005.30"
$ws.Range("E2").Value = "  +1.15%  "

$ws.Range("D3").Value = "2.063.08"
$ws.Range("E3").Value = "  -1.39%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'249.73"
$ws.Range("E5").Value = "  -0.97%  "

$ws.Range("E6").Value = "  +2.25%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'55.09"
$ws.Range("E8").Value = "  +12.35%  "

$ws.Range("D9").Value = "'60.86"
$ws.Range("E9").Value = "  +1.00%  "

$ws.Range("D10").Value = "'0.382"
$ws.Range("E10").Value = "  +1.85%  "

$ws.Range("D11").Value = "'0.0804"
$ws.Range("E11").Value = "  +7.86%  "

$ws.Range("E12").Value = "  +5.96%  "

$ws.Range("D13").Value = "'15.04"
$ws.Range("E13").Value = "  +1.76%  "

$ws.Range("D14").Value = "2.365.47"
$ws.Range("E14").Value = "  -1.41%  "

$ws.Range("E15").Value = "  -1.90%  "

$ws.Range("E16").Value = "  +3.91%  "

$ws.Range("D17").Value = "2.063.07"
$ws.Range("E17").Value = "  -1.23%  "

$ws.Range("D18").Value = "36.962.33"
$ws.Range("E18").Value = "  +1.17%  "

$ws.Range("D19").Value = "0.0₃0950"
$ws.Range("E19").Value = "  +13.93%  "

$ws.Range("D20").Value = "'73.47"
$ws.Range("E20").Value = "  +0.61%  "

$ws.Range("D21").Value = "'14.19"
$ws.Range("E21").Value = "  +6.98%  "

$ws.Range("E22").Value = "  +2.30%  "

$ws.Range("D23").Value = "'237.66"
$ws.Range("E23").Value = "  -0.96%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("E25").Value = "  -3.57%  "

$ws.Range("D26").Value = "'174.84"
$ws.Range("E26").Value = "  +2.19%  "

$ws.Range("E27").Value = "  -1.51%  "

$ws.Range("D28").Value = "'20.15"
$ws.Range("E28").Value = "  -4.64%  "

$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("E30").Value = "  +2.16%  "

$ws.Range("D31").Value = "'4.60"
$ws.Range("E31").Value = "  +2.39%  "

$ws.Range("E32").Value = "  +8.39%  "

$ws.Range("D33").Value = "'0.0628"
$ws.Range("E33").Value = "  +1.88%  "

$ws.Range("E34").Value = "  +7.42%  "

$ws.Range("E35").Value = "  -0.90%  "

$ws.Range("E36").Value = "  -0.22%  "

$ws.Range("E37").Value = "  -5.64%  "

$ws.Range("E38").Value = "  -4.78%  "

$ws.Range("E39").Value = "  +0.10%  "

$ws.Range("E40").Value = "  +24.22%  "

$ws.Range("D41").Value = "'17.84"
$ws.Range("E41").Value = "  +8.15%  "

$ws.Range("E42").Value = "  +0.82%  "

$ws.Range("E43").Value = "  -1.78%  "

$ws.Range("E44").Value = "  -1.07%  "

$ws.Range("E45").Value = "  +0.44%  "

$ws.Range("D46").Value = "'4.10"
$ws.Range("E46").Value = "  +33.83%  "

$ws.Range("D47").Value = "'13.98"
$ws.Range("E47").Value = "  -50.87%  "

$ws.Range("E48").Value = "  +7.55%  "

$ws.Range("D49").Value = "'4.29"
$ws.Range("E49").Value = "  +11.24%  "

$ws.Range("D50").Value = "1.302.14"
$ws.Range("E50").Value = "  -2.54%  "

$ws.Range("D51").Value = "'2.92"
$ws.Range("E51").Value = "  +1.57%  "
